$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The reporting period column "6 ماهه منتهی به 1399/06" (column D) is dropped
# entirely: every later quarter's figures shift one column to the left, and a
# brand-new closing quarter ("12 ماهه منتهی به 1401/12", dated 1402-02-10) is
# appended as the new final column (M).
$ws.Columns("D").Delete()

# Column L's formatting (number formats / fills / borders / fonts) is what the
# new trailing column M should inherit, since in this sheet every data row
# uses one uniform style across all of its quarter columns.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)  # xlPasteFormats

# Header rows: period label / publish-date label for the newest quarter.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-10 (2)"

# Financial data for the new quarter.
$ws.Range("M11").Value = 31261
$ws.Range("M12").Value = -24885
$ws.Range("M13").Value = 6375
$ws.Range("M14").Value = -612
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 94
$ws.Range("M17").Value = 5858
$ws.Range("M18").Value = -38
$ws.Range("M19").Value = -36
$ws.Range("M20").Value = 5783
$ws.Range("M21").Value = -594
$ws.Range("M22").Value = 5189
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 5189
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 2692
$ws.Range("M27").Value = 0
